$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.928222119808197
$ws.Range("B1").Value = 0.8046717643737793
$ws.Range("C1").Value = 0.7049208879470825
$ws.Range("D1").Value = 0.7624741792678833
$ws.Range("E1").Value = 0.929070770740509
